# Fruta / hortaliza, semanal
# Weekly update: insert two new price records near row 974 (pushing the
# existing rows down by two) and append two more records at the end of
# the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new rows at 974, shifting the rest of the table down ---
$ws.Rows.Item(974).Insert()
$ws.Rows.Item(974).Insert()

# Shared "boilerplate" values repeated down the whole table.
$mercadoId   = 9
$mercado     = "Vega Central Mapocho de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100102
$producto    = "Cítricos"
$categoriaId = 100102005
$categoria   = "Naranja"

# --- New row 974: Lane Late, Cabildo ---
$ws.Cells.Item(974, 1).Value  = $mercadoId
$ws.Cells.Item(974, 2).Value  = $mercado
$ws.Cells.Item(974, 3).Value  = $region
$ws.Cells.Item(974, 4).Value  = 44931
$ws.Cells.Item(974, 5).Value  = $codreg
$ws.Cells.Item(974, 6).Value  = $tipo
$ws.Cells.Item(974, 7).Value  = $productoId
$ws.Cells.Item(974, 8).Value  = $producto
$ws.Cells.Item(974, 9).Value  = $categoriaId
$ws.Cells.Item(974, 10).Value = $categoria
$ws.Cells.Item(974, 11).Value = "Lane Late"
$ws.Cells.Item(974, 12).Value = "Primera"
$ws.Cells.Item(974, 13).Value = 350
$ws.Cells.Item(974, 14).Value = 9000
$ws.Cells.Item(974, 15).Value = 9000
$ws.Cells.Item(974, 16).Value = 9000
$ws.Cells.Item(974, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(974, 18).Value = "Cabildo"
$ws.Cells.Item(974, 19).Value = 600
$ws.Cells.Item(974, 20).Value = 15

# --- New row 975: Navel Late, Provincia de Melipilla ---
$ws.Cells.Item(975, 1).Value  = $mercadoId
$ws.Cells.Item(975, 2).Value  = $mercado
$ws.Cells.Item(975, 3).Value  = $region
$ws.Cells.Item(975, 4).Value  = 44931
$ws.Cells.Item(975, 5).Value  = $codreg
$ws.Cells.Item(975, 6).Value  = $tipo
$ws.Cells.Item(975, 7).Value  = $productoId
$ws.Cells.Item(975, 8).Value  = $producto
$ws.Cells.Item(975, 9).Value  = $categoriaId
$ws.Cells.Item(975, 10).Value = $categoria
$ws.Cells.Item(975, 11).Value = "Navel Late"
$ws.Cells.Item(975, 12).Value = "Primera"
$ws.Cells.Item(975, 13).Value = 400
$ws.Cells.Item(975, 14).Value = 11000
$ws.Cells.Item(975, 15).Value = 11000
$ws.Cells.Item(975, 16).Value = 11000
$ws.Cells.Item(975, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(975, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(975, 19).Value = 611
$ws.Cells.Item(975, 20).Value = 18

# --- 2. Append two more records at the end of the table ---

# New row 1051: New Hall, Primera, Región de O'Higgins
$ws.Cells.Item(1051, 1).Value  = $mercadoId
$ws.Cells.Item(1051, 2).Value  = $mercado
$ws.Cells.Item(1051, 3).Value  = $region
$ws.Cells.Item(1051, 4).Value  = 44412
$ws.Cells.Item(1051, 5).Value  = $codreg
$ws.Cells.Item(1051, 6).Value  = $tipo
$ws.Cells.Item(1051, 7).Value  = $productoId
$ws.Cells.Item(1051, 8).Value  = $producto
$ws.Cells.Item(1051, 9).Value  = $categoriaId
$ws.Cells.Item(1051, 10).Value = $categoria
$ws.Cells.Item(1051, 11).Value = "New Hall"
$ws.Cells.Item(1051, 12).Value = "Primera"
$ws.Cells.Item(1051, 13).Value = 220
$ws.Cells.Item(1051, 14).Value = 4000
$ws.Cells.Item(1051, 15).Value = 4000
$ws.Cells.Item(1051, 16).Value = 4000
$ws.Cells.Item(1051, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(1051, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1051, 19).Value = 267
$ws.Cells.Item(1051, 20).Value = 15

# New row 1052: New Hall, Segunda, Región de O'Higgins
$ws.Cells.Item(1052, 1).Value  = $mercadoId
$ws.Cells.Item(1052, 2).Value  = $mercado
$ws.Cells.Item(1052, 3).Value  = $region
$ws.Cells.Item(1052, 4).Value  = 44412
$ws.Cells.Item(1052, 5).Value  = $codreg
$ws.Cells.Item(1052, 6).Value  = $tipo
$ws.Cells.Item(1052, 7).Value  = $productoId
$ws.Cells.Item(1052, 8).Value  = $producto
$ws.Cells.Item(1052, 9).Value  = $categoriaId
$ws.Cells.Item(1052, 10).Value = $categoria
$ws.Cells.Item(1052, 11).Value = "New Hall"
$ws.Cells.Item(1052, 12).Value = "Segunda"
$ws.Cells.Item(1052, 13).Value = 240
$ws.Cells.Item(1052, 14).Value = 3500
$ws.Cells.Item(1052, 15).Value = 3500
$ws.Cells.Item(1052, 16).Value = 3500
$ws.Cells.Item(1052, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(1052, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1052, 19).Value = 233
$ws.Cells.Item(1052, 20).Value = 15

Write-Output "Rows 974-975 inserted and rows 1051-1052 appended."
